# Daily scrape update - 2025-09-05 03:06:19 UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Replace the data rows (2-10) with the freshly scraped records, and
#    delete the now-stale rows 11-14 (the sheet shrinks from 14 to 9 data
#    rows + header).
# ---------------------------------------------------------------------------

$data = @(
    @("1327365", "https://aiesec.org/opportunity/global-talent/1327365", "[DSC] Finance Data Analyst Intern", "Fritz-Erler-Straße 5, 53113 Bonn, Germany", "Yes", "10 applicants", "6 - 18 Months", "DHL Group"),
    @("1327293", "https://aiesec.org/opportunity/global-talent/1327293", "Business Executive", "Hong Kong", "No", "9 applicants", "6 - 18 Months", "ASA Building Materials (HK) Limited"),
    @("1327242", "https://aiesec.org/opportunity/global-talent/1327242", "Sales Manager", "Cairo, Cairo Governorate, Egypt", "No", "0 applicants", "9 - 12 Weeks", "MZ creatives"),
    @("1327239", "https://aiesec.org/opportunity/global-talent/1327239", "Marketing Specialist", "Maadi as Sarayat Al Gharbeyah, Maadi, Cairo Governorate, Egypt", "No", "2 applicants", "3 - 6 Months", "Karcel"),
    @("1326449", "https://aiesec.org/opportunity/global-talent/1326449", "Treasury Analyst", "Panamá, Provincia de Panamá, Panamá", "No", "65 applicants", "6 - 18 Months", "Samsung Electronics Latinoamérica (Zona Libre) S.A (SELA)"),
    @("1324910", "https://aiesec.org/opportunity/global-talent/1324910", "Graphic designer", "Cairo, Cairo Governorate, Egypt", "No", "8 applicants", "3 - 6 Months", "Transition Agency"),
    @("1323761", "https://aiesec.org/opportunity/global-talent/1323761", "Sales representative", "Mansoura, Mansoura Qism 2, El Mansoura, Dakahlia Governorate, Egypt", "No", "4 applicants", "9 - 12 Weeks", "Fekretk"),
    @("1317293", "https://aiesec.org/opportunity/global-talent/1317293", "Guest Relations Executive", "Weligama, Sri Lanka", "No", "68 applicants", "9 - 12 Weeks", "Steradian Capital Investments"),
    @("1316788", "https://aiesec.org/opportunity/global-talent/1316788", "Travel Coordinator", "Mexico City, CDMX, Mexico", "No", "98 applicants", "6 - 18 Months", "Ikan Experience")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $rec = $data[$i]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
    $ws.Cells.Item($row, 8).Value = $rec[7]
}

# Highlight the new "Premium = Yes" opportunity (row 2, column E) in yellow.
$ws.Cells.Item(2, 5).Interior.Color = 65535

# Drop the 4 oldest listings - the sheet now only spans down to row 10.
$ws.Range("A11:H14").EntireRow.Delete()

# ---------------------------------------------------------------------------
# 2. Re-tune a few column widths to fit the new content.
# ---------------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 36
$ws.Columns.Item(4).ColumnWidth = 70
$ws.Columns.Item(6).ColumnWidth = 16
$ws.Columns.Item(8).ColumnWidth = 60
